$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row (54): raw/clean SSA COVID data for 2020-07-23.
# Column A holds a date-like label ("2020-07-23") that must stay a
# plain text value, matching every other row in the column, instead of
# being auto-converted by Excel into a date serial number. Forcing the
# cell to Text format before writing the value prevents that
# auto-conversion; resetting the cell style back to "Normal" afterwards
# keeps the cell's style index the same as all of its neighbours (no
# explicit s="" attribute), matching the original authoring style.
$ws.Range("A54").NumberFormat = "@"
$ws.Range("A54").Value = "2020-07-23"
$ws.Range("A54").Style = "Normal"

$ws.Range("B54").Value = 370712
$ws.Range("C54").Value = 419349
$ws.Range("D54").Value = 89547
$ws.Range("E54").Value = 41908
$ws.Range("F54").Value = 28.12
